$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.456.55"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "'1.588.34"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +1.07%  "
$ws.Range("D5").Value = "'213.23"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'0.493"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("D8").Value = "'24.51"
$ws.Range("E8").Value = "  +8.17%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'0.0600"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "'0.0886"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").Value = "'1.815.34"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "'1.592.72"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "'3.75"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "'28.432.27"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").Value = "'63.13"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "'230.12"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'7.49"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'0.0₃0707"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D22").Value = "'4.07"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "'9.36"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "'1.95"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'151.90"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "'15.22"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").Value = "'6.57"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'1.400.27"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -9.05%  "
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "'2.54"
$ws.Range("E39").Value = "  +8.52%  "
$ws.Range("D40").Value = "'0.541"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'0.811"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").Value = "'1.87"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "'5.58"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").Value = "'64.13"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "'1.725.71"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").Value = "'2.13"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "'87.18"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  +16.06%  "
